$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("od1")

# Fix the date strings so they read MM/DD/YYYY instead of DD/MM/YYYY
$ws.Range("A21").Value = "06/01/2018"
$ws.Range("A22").Value = "06/19/2018"
$ws.Range("A23").Value = "06/29/2018"

# Fix missing space between "29th" and "June" in the text body
$ws.Range("D23").Value = "<div class=""timeline-date"">29th June 2018</div> Almamet India signed an agreement for operation and maintenance including supply of reagents and consumables for HMDS facility with Jindal Steel and Power (JSPL) for its Angul works. `nThis makes Almamet unique as we supply engineering & equipment for HMDS and thereafter also take care of complete operation and maintenance. "
